$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-6 (years 2005年-2009年), so that the rows for
# 2010年-2015年 shift up from rows 7-12 to become rows 2-7.
$ws.Range("A2:G6").EntireRow.Delete()
